$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.114.85"
$ws.Range("E2").Value = "  -5.08%  "
$ws.Range("D3").Value = "3.312.87"
$ws.Range("E3").Value = "  -5.38%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "563.23"
$ws.Range("E5").Value = "  -4.49%  "
$ws.Range("D6").Value = "127.71"
$ws.Range("E6").Value = "  -4.78%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.312.62"
$ws.Range("E8").Value = "  -5.37%  "
$ws.Range("D9").Value = "0.478"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("E10").Value = "  -4.99%  "
$ws.Range("E11").Value = "  -5.21%  "
$ws.Range("E12").Value = "  -3.48%  "
$ws.Range("D13").Value = "3.882.56"
$ws.Range("E13").Value = "  -5.22%  "
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "3.315.67"
$ws.Range("E15").Value = "  -5.35%  "
$ws.Range("E16").Value = "  -6.71%  "
$ws.Range("D18").Value = "61.200.02"
$ws.Range("E18").Value = "  -4.87%  "
$ws.Range("D19").Value = "13.42"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value = "5.63"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("D21").Value = "8.98"
$ws.Range("E21").Value = "  -10.49%  "
$ws.Range("D22").Value = "352.36"
$ws.Range("E22").Value = "  -8.89%  "
$ws.Range("D23").Value = "0.554"
$ws.Range("E23").Value = "  -4.37%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "3.446.41"
$ws.Range("E25").Value = "  -5.36%  "
$ws.Range("D26").Value = "69.19"
$ws.Range("E26").Value = "  -6.99%  "
$ws.Range("E27").Value = "  -7.90%  "
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").Value = "7.14"
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("D30").Value = "7.85"
$ws.Range("E30").Value = "  -3.70%  "
$ws.Range("D31").Value = "1.40"
$ws.Range("E31").Value = "  -7.59%  "
$ws.Range("E32").Value = "  -6.96%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "0.149"
$ws.Range("E34").Value = "  -3.93%  "
$ws.Range("D35").Value = "3.340.09"
$ws.Range("E35").Value = "  -5.41%  "
$ws.Range("D36").Value = "22.53"
$ws.Range("E36").Value = "  -3.32%  "
$ws.Range("E37").Value = "  -3.79%  "
$ws.Range("D38").Value = "6.76"
$ws.Range("E38").Value = "  -2.14%  "
$ws.Range("D39").Value = "160.58"
$ws.Range("E39").Value = "  -2.69%  "
$ws.Range("D40").Value = "1.47"
$ws.Range("E40").Value = "  -4.57%  "
$ws.Range("D41").Value = "0.0756"
$ws.Range("E41").Value = "  -4.08%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "41.01"
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("D44").Value = "4.33"
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("D45").Value = "0.743"
$ws.Range("E45").Value = "  -8.01%  "
$ws.Range("D46").Value = "1.11"
$ws.Range("E46").Value = "  -6.29%  "
$ws.Range("E47").Value = "  -5.69%  "
$ws.Range("D48").Value = "22.25"
$ws.Range("E48").Value = "  -8.32%  "
$ws.Range("D49").Value = "6.68"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("D50").Value = "0.862"
$ws.Range("E50").Value = "  -6.45%  "
$ws.Range("D51").Value = "20.87"
$ws.Range("E51").Value = "  -1.08%  "